$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 111079.468651622
$ws.Range("E3").Value = -0.01116071249361237
$ws.Range("F3").Value = 0.188395706625901
$ws.Range("G3").Value = -1.222608724358497
$ws.Range("H3").Value = 10.5472320775076

# Row 5
$ws.Range("D5").Value = 112446.9664152832
$ws.Range("E5").Value = -0.0194313902090042
$ws.Range("F5").Value = 0.2287055213250087
$ws.Range("G5").Value = -1.313600381246777
$ws.Range("H5").Value = 10.42141036965161

# Row 7
$ws.Range("D7").Value = 114247.3370773786
$ws.Range("E7").Value = -0.02394101936143448
$ws.Range("F7").Value = 0.2407268286665641
$ws.Range("G7").Value = -1.139285595086238
$ws.Range("H7").Value = 8.222445247610166

# Row 8
$ws.Range("D8").Value = 114177.2107535922
$ws.Range("E8").Value = -0.03435496483924053
$ws.Range("F8").Value = 0.2019042775766292
$ws.Range("G8").Value = -1.251704475396882
$ws.Range("H8").Value = 8.815998366658114

# Row 9
$ws.Range("D9").Value = 115699.798632709
$ws.Range("E9").Value = -0.06015601302372368
$ws.Range("F9").Value = 0.3045045089566139
$ws.Range("G9").Value = -1.538128228484301
$ws.Range("H9").Value = 9.746143123177173

# Row 10
$ws.Range("D10").Value = 117188.6288699198
$ws.Range("E10").Value = -0.1029460749259843
$ws.Range("F10").Value = 0.4226795384216328
$ws.Range("G10").Value = -1.885997630692881
$ws.Range("H10").Value = 9.500789903146227

# Row 11
$ws.Range("D11").Value = 119043.3840067312
$ws.Range("E11").Value = -0.1729122259264516
$ws.Range("F11").Value = 0.708835519173359
$ws.Range("G11").Value = -2.463814158598304
$ws.Range("H11").Value = 11.47676327413394

# Row 13
$ws.Range("D13").Value = 110252.75002931
$ws.Range("E13").Value = 0.006930772844332807
$ws.Range("F13").Value = 0.1706857216703247
$ws.Range("G13").Value = -0.7135665547298811
$ws.Range("H13").Value = 6.170919588763635
